$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: append a new record under the existing table (A1:H3 -> A1:H4).
# Leading "'" forces text storage (empty string for A4, "2222" kept
# as text rather than becoming a number for C4), matching the other
# rows' t="str" cell type.
$ws.Range("A4").Value = "'"
$ws.Range("B4").Value = "حسن "
$ws.Range("C4").Value = "'2222"
$ws.Range("D4").Value = "ايتا"
$ws.Range("E4").Value = "الرحلة 2"
$ws.Range("F4").Value = "C3"
$ws.Range("G4").Value = "NRC"
$ws.Range("H4").Value = "٠٢‏/٠٥‏/٢٠٢٥ ٠١:٥٢:٥١ م"

# The quote-prefix entry above stamps a transient "quotePrefix" style on
# A4/C4; reset the whole new row back to the default (unstyled) look so
# it matches rows 1-3, which carry no style attribute at all.
$ws.Range("A4:H4").Style = "Normal"
